# Modified LoginPage POM: added an "InvalidLogin" worksheet (with its own
# UserName/Password header row and a Bhanu/Damager bad-credentials row),
# changed the selection on the original "ValidLogin" sheet, and made the
# new sheet the active tab.

$wb = $excel.ActiveWorkbook

# --- ValidLogin: update the remembered selection (was C2, now A1:B2) ---
$validLogin = $wb.Worksheets.Item("ValidLogin")
$validLogin.Range("A1:B2").Select() | Out-Null

# --- Add the new InvalidLogin sheet after ValidLogin ---
$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

# Header row matches ValidLogin's, data row holds bad credentials.
$invalidLogin.Range("A1").Value = "UserName"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("A2").Value = "Bhanu"
$invalidLogin.Range("B2").Value = "Damager"

# Auto-size the columns to the new content.
$invalidLogin.Range("A:B").AutoFit() | Out-Null

# Zoom in on the new sheet and leave the selection on B3, then make it active.
$excel.ActiveWindow.Zoom = 250
$invalidLogin.Range("B3").Select() | Out-Null
$invalidLogin.Activate() | Out-Null
